# Weekly update: insert 3 new rows (new week's Palta data) above the
# existing block starting at row 416, pushing the previous rows down by
# three (416->419 ... 436->439), matching a "new week pushed to top,
# older weeks shift down" ingestion pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 416; Excel shifts 416:436 -> 419:439.
$ws.Rows("416:418").Insert()

# New week (2022-05-25, serial 44706) data for Comercializadora del Agro
# de Limarí - Palta - Hass, categories Especial / Primera / Segunda.

# Row 416 - Hass Especial
$ws.Range("A416").Value = 2
$ws.Range("B416").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C416").Value = "Coquimbo"
$ws.Range("D416").Value = "2022-05-25"
$ws.Range("E416").Value = 4
$ws.Range("F416").Value = "Fruta"
$ws.Range("G416").Value = 100106
$ws.Range("H416").Value = "Oleaginosos"
$ws.Range("I416").Value = 100106002
$ws.Range("J416").Value = "Palta"
$ws.Range("K416").Value = "Hass"
$ws.Range("L416").Value = "Especial"
$ws.Range("M416").Value = 400
$ws.Range("N416").Value = 2900
$ws.Range("O416").Value = 3000
$ws.Range("P416").Value = 2950
$ws.Range("Q416").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R416").Value = "Provincia de Limarí"
$ws.Range("S416").Value = 2950
$ws.Range("T416").Value = 1

# Row 417 - Hass Primera
$ws.Range("A417").Value = 2
$ws.Range("B417").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C417").Value = "Coquimbo"
$ws.Range("D417").Value = "2022-05-25"
$ws.Range("E417").Value = 4
$ws.Range("F417").Value = "Fruta"
$ws.Range("G417").Value = 100106
$ws.Range("H417").Value = "Oleaginosos"
$ws.Range("I417").Value = 100106002
$ws.Range("J417").Value = "Palta"
$ws.Range("K417").Value = "Hass"
$ws.Range("L417").Value = "Primera"
$ws.Range("M417").Value = 400
$ws.Range("N417").Value = 2700
$ws.Range("O417").Value = 2800
$ws.Range("P417").Value = 2750
$ws.Range("Q417").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R417").Value = "Provincia de Limarí"
$ws.Range("S417").Value = 2750
$ws.Range("T417").Value = 1

# Row 418 - Hass Segunda
$ws.Range("A418").Value = 2
$ws.Range("B418").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C418").Value = "Coquimbo"
$ws.Range("D418").Value = "2022-05-25"
$ws.Range("E418").Value = 4
$ws.Range("F418").Value = "Fruta"
$ws.Range("G418").Value = 100106
$ws.Range("H418").Value = "Oleaginosos"
$ws.Range("I418").Value = 100106002
$ws.Range("J418").Value = "Palta"
$ws.Range("K418").Value = "Hass"
$ws.Range("L418").Value = "Segunda"
$ws.Range("M418").Value = 300
$ws.Range("N418").Value = 2500
$ws.Range("O418").Value = 2600
$ws.Range("P418").Value = 2550
$ws.Range("Q418").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R418").Value = "Provincia de Limarí"
$ws.Range("S418").Value = 2550
$ws.Range("T418").Value = 1
